$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new log rows (79-81) below the existing data (which ends at row 78)
$ws.Cells.Item(79, 1).Value = "2023-12-07 18:17:52"
$ws.Cells.Item(79, 2).Value = 0.002

$ws.Cells.Item(80, 1).Value = "2023-12-07 18:18:47"
$ws.Cells.Item(80, 2).Value = 0.003

$ws.Cells.Item(81, 1).Value = "2023-12-07 18:19:05"
$ws.Cells.Item(81, 2).Value = 0.0006000000000000001
